$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 544.25
$ws.Range("J29").Value = 202
$ws.Range("L29").Value = 606
$ws.Range("N29").Value = -1168

$ws.Range("H38").Value = 358.41666
$ws.Range("I38").Value = 48.3125
$ws.Range("J38").Value = 978.625
$ws.Range("K38").Value = 144.9375
$ws.Range("L38").Value = 2935.875
$ws.Range("M38").Value = 227.0625
$ws.Range("N38").Value = -3679.875

$ws.Range("H40").Value = 1761.5454
$ws.Range("I40").Value = 1646.8695
$ws.Range("J40").Value = 2025.3
$ws.Range("K40").Value = 1646.8695
$ws.Range("L40").Value = 2025.3
$ws.Range("M40").Value = -1471.8695
$ws.Range("N40").Value = -2375.3

$ws.Range("H43").Value = 6125
$ws.Range("I43").Value = 6000
$ws.Range("J43").Value = 6250
$ws.Range("K43").Value = 6000
$ws.Range("L43").Value = 6250
$ws.Range("M43").Value = -5931
$ws.Range("N43").Value = -6388

$ws.Range("H58").Value = 998.5714
$ws.Range("I58").Value = 376.5
$ws.Range("J58").Value = 1247.4
$ws.Range("K58").Value = 1129.5
$ws.Range("L58").Value = 3742.2
$ws.Range("M58").Value = -979.5
$ws.Range("N58").Value = -4042.2

$ws.Range("H86").Value = 2055.6428
$ws.Range("I86").Value = 1997.7778
$ws.Range("J86").Value = 2159.8
$ws.Range("K86").Value = 1997.7778
$ws.Range("L86").Value = 2159.8
$ws.Range("M86").Value = -874.7778000000001
$ws.Range("N86").Value = -4405.8

$ws.Range("H87").Value = 35195
$ws.Range("J87").Value = 35195
$ws.Range("L87").Value = 35195
$ws.Range("N87").Value = -37691

$ws.Range("H89").Value = 2055.6428
$ws.Range("I89").Value = 1997.7778
$ws.Range("J89").Value = 2159.8
$ws.Range("K89").Value = 9988.889000000001
$ws.Range("L89").Value = 10799
$ws.Range("M89").Value = -4372.889000000001
$ws.Range("N89").Value = -22031

$ws.Range("H90").Value = 35195
$ws.Range("J90").Value = 35195
$ws.Range("L90").Value = 105585
$ws.Range("N90").Value = -118065

$ws.Range("H113").Value = 11289.2
$ws.Range("I113").Value = 41520.8
$ws.Range("J113").Value = 3731.3
$ws.Range("K113").Value = 41520.8
$ws.Range("L113").Value = 3731.3
$ws.Range("M113").Value = -38266.8
$ws.Range("N113").Value = -10239.3

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 27779032
$ws.Range("I122").Value = 35715180
$ws.Range("J122").Value = 2512.375
$ws.Range("K122").Value = 107145540
$ws.Range("L122").Value = 7537.125
$ws.Range("M122").Value = -107143090
$ws.Range("N122").Value = -12437.125

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 60492.9
$ws.Range("I80").Value = 275500
$ws.Range("J80").Value = 6741.125
$ws.Range("K80").Value = 275500
$ws.Range("L80").Value = 6741.125
$ws.Range("M80").Value = -274502
$ws.Range("N80").Value = -8737.125

$ws.Range("H83").Value = 60492.9
$ws.Range("I83").Value = 275500
$ws.Range("J83").Value = 6741.125
$ws.Range("K83").Value = 1377500
$ws.Range("L83").Value = 33705.625
$ws.Range("M83").Value = -1372508
$ws.Range("N83").Value = -43689.625

$ws.Range("H86").Value = 1488.56
$ws.Range("I86").Value = 1353.7059
$ws.Range("J86").Value = 1775.125
$ws.Range("K86").Value = 1353.7059
$ws.Range("L86").Value = 1775.125
$ws.Range("M86").Value = -230.7058999999999
$ws.Range("N86").Value = -4021.125

$ws.Range("H89").Value = 1488.56
$ws.Range("I89").Value = 1353.7059
$ws.Range("J89").Value = 1775.125
$ws.Range("K89").Value = 6768.5295
$ws.Range("L89").Value = 8875.625
$ws.Range("M89").Value = -1152.5295
$ws.Range("N89").Value = -20107.625

$ws.Range("H107").Value = 41668156
$ws.Range("I107").Value = 52633050
$ws.Range("J107").Value = 1578
$ws.Range("K107").Value = 52633050
$ws.Range("L107").Value = 1578
$ws.Range("M107").Value = -52631130
$ws.Range("N107").Value = -5418

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 8662.227999999999
$ws.Range("I31").Value = 1407.5
$ws.Range("J31").Value = 21358
$ws.Range("K31").Value = 1407.5
$ws.Range("L31").Value = 21358
$ws.Range("M31").Value = -1112.5
$ws.Range("N31").Value = -21948

$ws.Range("H34").Value = 8662.227999999999
$ws.Range("I34").Value = 1407.5
$ws.Range("J34").Value = 21358
$ws.Range("K34").Value = 1407.5
$ws.Range("L34").Value = 21358
$ws.Range("M34").Value = -1205.5
$ws.Range("N34").Value = -21762

$ws.Range("H122").Value = 1757.1765
$ws.Range("I122").Value = 1134.1538
$ws.Range("J122").Value = 3782
$ws.Range("K122").Value = 3402.4614
$ws.Range("L122").Value = 11346
$ws.Range("M122").Value = -952.4614000000001
$ws.Range("N122").Value = -16246

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 409.67856
$ws.Range("I5").Value = 311.625
$ws.Range("J5").Value = 540.4167
$ws.Range("K5").Value = 934.875
$ws.Range("L5").Value = 1621.2501
$ws.Range("M5").Value = -822.875
$ws.Range("N5").Value = -1845.2501

$ws.Range("H34").Value = 720.05884
$ws.Range("I34").Value = 139.28572
$ws.Range("J34").Value = 1126.6
$ws.Range("K34").Value = 417.85716
$ws.Range("L34").Value = 3379.8
$ws.Range("M34").Value = -333.85716
$ws.Range("N34").Value = -3547.8

$ws.Range("H39").Value = 1820.1
$ws.Range("I39").Value = 350
$ws.Range("J39").Value = 1897.4736
$ws.Range("K39").Value = 1050
$ws.Range("L39").Value = 5692.4208
$ws.Range("M39").Value = -756
$ws.Range("N39").Value = -6280.4208

$ws.Range("H55").Value = 920.2353000000001
$ws.Range("I55").Value = 592
$ws.Range("J55").Value = 964
$ws.Range("K55").Value = 1776
$ws.Range("L55").Value = 2892
$ws.Range("M55").Value = -1599
$ws.Range("N55").Value = -3246

$ws.Range("H60").Value = 1287.1666
$ws.Range("I60").Value = 1262.8889
$ws.Range("J60").Value = 1360
$ws.Range("K60").Value = 3788.6667
$ws.Range("L60").Value = 4080
$ws.Range("M60").Value = -3537.6667
$ws.Range("N60").Value = -4582

$ws.Range("H135").Value = 409.67856
$ws.Range("I135").Value = 311.625
$ws.Range("J135").Value = 540.4167
$ws.Range("K135").Value = 2804.625
$ws.Range("L135").Value = 4863.7503
$ws.Range("M135").Value = -269.625
$ws.Range("N135").Value = -9933.7503

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 7333.3335
$ws.Range("I43").Value = 4000
$ws.Range("J43").Value = 7750
$ws.Range("K43").Value = 4000
$ws.Range("L43").Value = 7750
$ws.Range("M43").Value = -3849
$ws.Range("N43").Value = -8052

$ws.Range("H46").Value = 14666.667
$ws.Range("J46").Value = 18000
$ws.Range("L46").Value = 18000
$ws.Range("N46").Value = -18312

$ws.Range("H57").Value = 6061
$ws.Range("J57").Value = 6061
$ws.Range("L57").Value = 6061
$ws.Range("N57").Value = -7701

$ws.Range("H80").Value = 1889.9111
$ws.Range("I80").Value = 1821.3077
$ws.Range("J80").Value = 1983.7894
$ws.Range("K80").Value = 1821.3077
$ws.Range("L80").Value = 1983.7894
$ws.Range("M80").Value = -823.3077000000001
$ws.Range("N80").Value = -3979.7894

$ws.Range("H83").Value = 1889.9111
$ws.Range("I83").Value = 1821.3077
$ws.Range("J83").Value = 1983.7894
$ws.Range("K83").Value = 9106.538500000001
$ws.Range("L83").Value = 9918.947
$ws.Range("M83").Value = -4114.538500000001
$ws.Range("N83").Value = -19902.947

$ws.Range("H113").Value = 1302.2413
$ws.Range("I113").Value = 1069.8096
$ws.Range("J113").Value = 1912.375
$ws.Range("K113").Value = 1069.8096
$ws.Range("L113").Value = 1912.375
$ws.Range("M113").Value = 1100.1904
$ws.Range("N113").Value = -6252.375

$ws.Range("H132").Value = 22753512
$ws.Range("I132").Value = 31282820
$ws.Range("J132").Value = 8696
$ws.Range("K132").Value = 93848460
$ws.Range("L132").Value = 26088
$ws.Range("M132").Value = -93845930
$ws.Range("N132").Value = -31148

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1651.5834
$ws.Range("I7").Value = 1482
$ws.Range("J7").Value = 2499.5
$ws.Range("K7").Value = 1482
$ws.Range("L7").Value = 2499.5
$ws.Range("M7").Value = -1370
$ws.Range("N7").Value = -2723.5

$ws.Range("H46").Value = 4057.4517
$ws.Range("I46").Value = 858.4167
$ws.Range("K46").Value = 858.4167
$ws.Range("M46").Value = -670.4167

$ws.Range("H61").Value = 1991.7858
$ws.Range("I61").Value = 1546.875
$ws.Range("J61").Value = 2585
$ws.Range("K61").Value = 1546.875
$ws.Range("L61").Value = 2585
$ws.Range("M61").Value = -1344.875
$ws.Range("N61").Value = -2989

$ws.Range("H113").Value = 1991.7858
$ws.Range("I113").Value = 1546.875
$ws.Range("J113").Value = 2585
$ws.Range("K113").Value = 1546.875
$ws.Range("L113").Value = 2585
$ws.Range("M113").Value = 623.125
$ws.Range("N113").Value = -6925

$ws.Range("H126").Value = 1651.5834
$ws.Range("I126").Value = 1482
$ws.Range("J126").Value = 2499.5
$ws.Range("K126").Value = 4446
$ws.Range("L126").Value = 7498.5
$ws.Range("M126").Value = -1976
$ws.Range("N126").Value = -12438.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H20").Value = 5468
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 5468
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 5468
$ws.Range("M20").ClearContents()
$ws.Range("N20").Value = -5948

$ws.Range("H126").Value = 50000596
$ws.Range("I126").Value = 14706229
$ws.Range("K126").Value = 44118687
$ws.Range("L126").Value = 375003360
$ws.Range("M126").Value = -44116217
$ws.Range("N126").Value = -375008300
